$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the rows containing "Huelva" and "Huesca" in column A and swap
# their province name and their "Casos activos" (column C) value, since
# the two rows' other figures (B, D, E) are identical.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$huelvaRow = -1
$huescaRow = -1

for ($r = 1; $r -le $lastRow; $r++) {
    $val = $ws.Cells.Item($r, 1).Value2
    if ($val -eq "Huelva") { $huelvaRow = $r }
    if ($val -eq "Huesca") { $huescaRow = $r }
}

if ($huelvaRow -gt 0 -and $huescaRow -gt 0) {
    $cHuelva = $ws.Cells.Item($huelvaRow, 3).Value2
    $cHuesca = $ws.Cells.Item($huescaRow, 3).Value2

    $ws.Cells.Item($huelvaRow, 1).Value2 = "Huesca"
    $ws.Cells.Item($huescaRow, 1).Value2 = "Huelva"

    $ws.Cells.Item($huelvaRow, 3).Value2 = $cHuesca
    $ws.Cells.Item($huescaRow, 3).Value2 = $cHuelva
}

# Update the "Datos actualizados" timestamp cell (A1) from 04:16 to 04:46.
$a1 = $ws.Range("A1")
$a1.Value2 = $a1.Value2 -replace "04:16", "04:46"
